$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "Indicator"

$indicatorValues = @(1,2,3,1,1,1,2,2,2,2,2,3,1,2,2,2,3,3,1)
for ($i = 0; $i -lt $indicatorValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $indicatorValues[$i]
}

$ws.Range("D21").Select()
